# Add a "Quantite" numeric column right after "Ingrédients" (new column C),
# shifting the previous C:G columns to D:H, and populate it with the numeric
# quantities that were previously encoded only as text (e.g. "4u" -> 4,
# "200g" -> 200, "100g" -> 100) in what is now column D ("Quantité Ingrédients").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hyperlinks that currently live in column G (they will end up in column H
# once the new column is inserted). Record their target row/URL now because
# inserting a column does not relocate the <hyperlinks> entries automatically
# in this runtime.
$hyperlinkInfo = @(
    @{ Row = 2; Url = "https://www.marmiton.org/shopping/envie-d-un-four-encastrable-au-top-voici-les-6-meilleurs-fours-avec-pyrolyse-s4008147.html?utm_source=ustensiles-recettes" },
    @{ Row = 3; Url = "https://www.marmiton.org/shopping/4-saladiers-qui-nous-ont-fait-craquer-s3032049.html?utm_source=ustensiles-recettes" },
    @{ Row = 6; Url = "https://www.marmiton.org/shopping/envie-d-un-four-encastrable-au-top-voici-les-6-meilleurs-fours-avec-pyrolyse-s4008147.html?utm_source=ustensiles-recettes" },
    @{ Row = 7; Url = "https://www.marmiton.org/shopping/envie-d-un-four-encastrable-au-top-voici-les-6-meilleurs-fours-avec-pyrolyse-s4008147.html?utm_source=ustensiles-recettes" }
)

# Drop the existing hyperlinks; they get recreated after the shift below.
$ws.Hyperlinks.Delete()

# Insert the new column before the old "Quantité Ingrédients" column (C),
# pushing C:G one column to the right, to D:H.
$ws.Columns("C:C").Insert()

# Header + numeric quantities for the freshly inserted column.
$ws.Range("C1").Value = "Quantite"
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 200
$ws.Range("C4").Value = 100
$ws.Range("C5").Value = 100
$ws.Range("C6").Value = 200

# Give the new column the same display width as column B.
$ws.Range("C1:C7").ColumnWidth = $ws.Range("B1").ColumnWidth

# Recreate the hyperlinks on column H (shifted right from the former column
# G). Hyperlinks.Add replaces the cell's content with the link text and
# re-applies the default "Hyperlink" look, so immediately afterwards we
# restore the original instruction sentence and copy back the plain
# (non-underlined / automatic-color) formatting the cell had before, taken
# from row 4 which was never touched by a hyperlink.
$ws.Range("H4").Copy() | Out-Null

foreach ($info in $hyperlinkInfo) {
    $cell = $ws.Cells.Item($info.Row, 8)
    $originalText = $cell.Value2

    $ws.Hyperlinks.Add($cell, $info.Url, "", "", $info.Url) | Out-Null

    $cell.Value = $originalText
    $cell.PasteSpecial(-4122) | Out-Null
}

# Restore the active selection, which moved from B6 to C6.
$ws.Range("C6").Select() | Out-Null
